$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.070935
$ws.Cells.Item(2, 8).Value = 6.212805
$ws.Cells.Item(2, 9).Value = 0.09632749399019591
$ws.Cells.Item(2, 10).Value = 0.09632749399019594
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.060838
$ws.Cells.Item(2, 14).Value = 0.182514
$ws.Cells.Item(2, 15).Value = 0.001214238342664256
$ws.Cells.Item(2, 16).Value = 0.001214238342664256
$ws.Cells.Item(2, 17).Value = 0.12599154353
$ws.Cells.Item(2, 18).Value = 1.13392389177
$ws.Cells.Item(2, 19).Value = 0.0001169645366556565
$ws.Cells.Item(2, 20).Value = 0.0001169645366556566

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.070935
$ws.Cells.Item(3, 8).Value = 6.212805
$ws.Cells.Item(3, 9).Value = 0.09632749399019591
$ws.Cells.Item(3, 10).Value = 0.09632749399019594
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 50.042999
$ws.Cells.Item(3, 14).Value = 150.128997
$ws.Cells.Item(3, 15).Value = 0.9987857616573358
$ws.Cells.Item(3, 16).Value = 0.9987857616573358
$ws.Cells.Item(3, 17).Value = 103.635798134065
$ws.Cells.Item(3, 18).Value = 932.722183206585
$ws.Cells.Item(3, 19).Value = 0.09621052945354026
$ws.Cells.Item(3, 20).Value = 0.09621052945354029

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.209141666666667
$ws.Cells.Item(4, 8).Value = 9.627425
$ws.Cells.Item(4, 9).Value = 0.1492700517445119
$ws.Cells.Item(4, 10).Value = 0.1492700517445119
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.060838
$ws.Cells.Item(4, 14).Value = 0.182514
$ws.Cells.Item(4, 15).Value = 0.001214238342664256
$ws.Cells.Item(4, 16).Value = 0.001214238342664256
$ws.Cells.Item(4, 17).Value = 0.1952377607166667
$ws.Cells.Item(4, 18).Value = 1.75713984645
$ws.Cells.Item(4, 19).Value = 0.0001812494202396637
$ws.Cells.Item(4, 20).Value = 0.0001812494202396638

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.209141666666667
$ws.Cells.Item(5, 8).Value = 9.627425
$ws.Cells.Item(5, 9).Value = 0.1492700517445119
$ws.Cells.Item(5, 10).Value = 0.1492700517445119
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 50.042999
$ws.Cells.Item(5, 14).Value = 150.128997
$ws.Cells.Item(5, 15).Value = 0.9987857616573358
$ws.Cells.Item(5, 16).Value = 0.9987857616573358
$ws.Cells.Item(5, 17).Value = 160.5950732158583
$ws.Cells.Item(5, 18).Value = 1445.355658942725
$ws.Cells.Item(5, 19).Value = 0.1490888023242722
$ws.Cells.Item(5, 20).Value = 0.1490888023242722

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.134776333333333
$ws.Cells.Item(6, 8).Value = 3.404329
$ws.Cells.Item(6, 9).Value = 0.05278299919088877
$ws.Cells.Item(6, 10).Value = 0.05278299919088877
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.060838
$ws.Cells.Item(6, 14).Value = 0.182514
$ws.Cells.Item(6, 15).Value = 0.001214238342664256
$ws.Cells.Item(6, 16).Value = 0.001214238342664256
$ws.Cells.Item(6, 17).Value = 0.06903752256733334
$ws.Cells.Item(6, 18).Value = 0.6213377031060001
$ws.Cells.Item(6, 19).Value = 0.00006409114145839352
$ws.Cells.Item(6, 20).Value = 0.00006409114145839352

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.134776333333333
$ws.Cells.Item(7, 8).Value = 3.404329
$ws.Cells.Item(7, 9).Value = 0.05278299919088877
$ws.Cells.Item(7, 10).Value = 0.05278299919088877
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 50.042999
$ws.Cells.Item(7, 14).Value = 150.128997
$ws.Cells.Item(7, 15).Value = 0.9987857616573358
$ws.Cells.Item(7, 16).Value = 0.9987857616573358
$ws.Cells.Item(7, 17).Value = 56.78761091422367
$ws.Cells.Item(7, 18).Value = 511.088498228013
$ws.Cells.Item(7, 19).Value = 0.05271890804943038
$ws.Cells.Item(7, 20).Value = 0.05271890804943039

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.08404533333333
$ws.Cells.Item(8, 8).Value = 45.252136
$ws.Cells.Item(8, 9).Value = 0.7016194550744034
$ws.Cells.Item(8, 10).Value = 0.7016194550744034
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.060838
$ws.Cells.Item(8, 14).Value = 0.182514
$ws.Cells.Item(8, 15).Value = 0.001214238342664256
$ws.Cells.Item(8, 16).Value = 0.001214238342664256
$ws.Cells.Item(8, 17).Value = 0.9176831499893334
$ws.Cells.Item(8, 18).Value = 8.259148349904
$ws.Cells.Item(8, 19).Value = 0.0008519332443105417
$ws.Cells.Item(8, 20).Value = 0.0008519332443105417

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.08404533333333
$ws.Cells.Item(9, 8).Value = 45.252136
$ws.Cells.Item(9, 9).Value = 0.7016194550744034
$ws.Cells.Item(9, 10).Value = 0.7016194550744034
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 50.042999
$ws.Cells.Item(9, 14).Value = 150.128997
$ws.Cells.Item(9, 15).Value = 0.9987857616573358
$ws.Cells.Item(9, 16).Value = 0.9987857616573358
$ws.Cells.Item(9, 17).Value = 754.8508655319547
$ws.Cells.Item(9, 18).Value = 6793.657789787592
$ws.Cells.Item(9, 19).Value = 0.7007675218300928
$ws.Cells.Item(9, 20).Value = 0.7007675218300928
